$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 2.5
$ws.Range("F3").Value = "Terminé"

$ws.Range("E4").Value = 1
$ws.Range("F4").Value = "En cours"

$ws.Range("F12").Value = "Terminé"

$ws.Range("G5").Select()
